# Scheduled market-data refresh: update currentAveragePrice* / Leve price &
# profit columns (H:N) for the affected Leve rows on each Job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 342.65
$ws.Range("I15").Value = 342.65
$ws.Range("K15").Value = 1027.95
$ws.Range("M15").Value = -858.9499999999998

$ws.Range("H19").Value = 1074
$ws.Range("I19").Value = 900
$ws.Range("J19").Value = 1132
$ws.Range("K19").Value = 900
$ws.Range("L19").Value = 1132
$ws.Range("M19").Value = -725
$ws.Range("N19").Value = -1482

$ws.Range("H64").Value = 5201.6665
$ws.Range("J64").Value = 4830.7144
$ws.Range("L64").Value = 4830.7144
$ws.Range("N64").Value = -5326.7144

$ws.Range("H67").Value = 5201.6665
$ws.Range("J67").Value = 4830.7144
$ws.Range("L67").Value = 4830.7144
$ws.Range("N67").Value = -6546.7144

$ws.Range("H99").Value = 836.7692
$ws.Range("I99").Value = 520
$ws.Range("J99").Value = 1549.5
$ws.Range("K99").Value = 1560
$ws.Range("L99").Value = 4648.5
$ws.Range("M99").Value = -62
$ws.Range("N99").Value = -7644.5

$ws.Range("H132").Value = 964726.8
$ws.Range("I132").Value = 3809
$ws.Range("J132").Value = 4904490
$ws.Range("K132").Value = 11427
$ws.Range("L132").Value = 14713470
$ws.Range("M132").Value = -8897
$ws.Range("N132").Value = -14718530

$ws.Range("H138").Value = 3392384.2
$ws.Range("I138").Value = 1680.0435
$ws.Range("J138").Value = 5558667.5
$ws.Range("K138").Value = 5040.1305
$ws.Range("L138").Value = 16676002.5
$ws.Range("M138").Value = 99.86949999999979
$ws.Range("N138").Value = -16686282.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4703.37
$ws.Range("I32").Value = 4415.2
$ws.Range("K32").Value = 4415.2
$ws.Range("M32").Value = -4128.2

$ws.Range("H61").Value = 77078520
$ws.Range("I61").Value = 83418390
$ws.Range("K61").Value = 83418390
$ws.Range("M61").Value = -83418178

$ws.Range("H97").Value = 2842362.5
$ws.Range("I97").Value = 3907648
$ws.Range("J97").Value = 1601.6666
$ws.Range("K97").Value = 3907648
$ws.Range("L97").Value = 1601.6666
$ws.Range("M97").Value = -3907152
$ws.Range("N97").Value = -2593.6666

$ws.Range("H102").Value = 20411348
$ws.Range("I102").Value = 28574828
$ws.Range("J102").Value = 2650
$ws.Range("K102").Value = 28574828
$ws.Range("L102").Value = 2650
$ws.Range("M102").Value = -28573206
$ws.Range("N102").Value = -5894

$ws.Range("H132").Value = 7970930.5
$ws.Range("I132").Value = 8639089
$ws.Range("K132").Value = 25917267
$ws.Range("M132").Value = -25914737

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 65908.625
$ws.Range("J135").Value = 65908.625
$ws.Range("L135").Value = 65908.625
$ws.Range("N135").Value = -76048.625

$ws.Range("H136").Value = 77078520
$ws.Range("I136").Value = 83418390
$ws.Range("K136").Value = 250255170
$ws.Range("M136").Value = -250252620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 697
$ws.Range("I94").Value = 697
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 697
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -246
$ws.Range("N94").ClearContents()

$ws.Range("H107").Value = 3598.5
$ws.Range("I107").Value = 2873.7144
$ws.Range("J107").Value = 4769.3076
$ws.Range("K107").Value = 2873.7144
$ws.Range("L107").Value = 4769.3076
$ws.Range("M107").Value = -953.7143999999998
$ws.Range("N107").Value = -8609.3076

$ws.Range("H134").Value = 9528808
$ws.Range("I134").Value = 5144.853
$ws.Range("J134").Value = 333333340
$ws.Range("K134").Value = 15434.559
$ws.Range("L134").Value = 1000000020
$ws.Range("M134").Value = -12899.559
$ws.Range("N134").Value = -1000005090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 259934.72
$ws.Range("I31").Value = 55024.633
$ws.Range("J31").Value = 398980.84
$ws.Range("K31").Value = 55024.633
$ws.Range("L31").Value = 398980.84
$ws.Range("M31").Value = -54729.633
$ws.Range("N31").Value = -399570.84

$ws.Range("H34").Value = 259934.72
$ws.Range("I34").Value = 55024.633
$ws.Range("J34").Value = 398980.84
$ws.Range("K34").Value = 55024.633
$ws.Range("L34").Value = 398980.84
$ws.Range("M34").Value = -54822.633
$ws.Range("N34").Value = -399384.84

$ws.Range("H99").Value = 1449.5
$ws.Range("I99").Value = 1066
$ws.Range("J99").Value = 1833
$ws.Range("K99").Value = 1066
$ws.Range("L99").Value = 1833
$ws.Range("M99").Value = 432
$ws.Range("N99").Value = -4829

$ws.Range("H126").Value = 1449.5
$ws.Range("I126").Value = 1066
$ws.Range("J126").Value = 1833
$ws.Range("K126").Value = 3198
$ws.Range("L126").Value = 5499
$ws.Range("M126").Value = -728
$ws.Range("N126").Value = -10439

$ws.Range("H132").Value = 73786.42999999999
$ws.Range("I132").Value = 1573.1428
$ws.Range("J132").Value = 145999.72
$ws.Range("K132").Value = 4719.428400000001
$ws.Range("L132").Value = 437999.16
$ws.Range("M132").Value = -2189.428400000001
$ws.Range("N132").Value = -443059.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 851.63635
$ws.Range("J131").Value = 1045.6666
$ws.Range("L131").Value = 3136.9998
$ws.Range("N131").Value = -13216.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1972.6818
$ws.Range("I102").Value = 1960.6875
$ws.Range("J102").Value = 2004.6666
$ws.Range("K102").Value = 1960.6875
$ws.Range("L102").Value = 2004.6666
$ws.Range("M102").Value = -338.6875
$ws.Range("N102").Value = -5248.6666

$ws.Range("H132").Value = 47403.5
$ws.Range("I132").Value = 30862.383
$ws.Range("J132").Value = 103643.3
$ws.Range("K132").Value = 92587.149
$ws.Range("L132").Value = 310929.9
$ws.Range("M132").Value = -90057.149
$ws.Range("N132").Value = -315989.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83335170
$ws.Range("I7").Value = 83335170
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 83335170
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -83335058
$ws.Range("N7").ClearContents()

$ws.Range("H122").Value = 6636.6665
$ws.Range("I122").Value = 11915
$ws.Range("J122").Value = 3997.5
$ws.Range("K122").Value = 35745
$ws.Range("L122").Value = 11992.5
$ws.Range("M122").Value = -33295
$ws.Range("N122").Value = -16892.5

$ws.Range("H126").Value = 83335170
$ws.Range("I126").Value = 83335170
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 250005510
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -250003040
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1234
$ws.Range("I122").Value = 880.8
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2642.4
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -192.3999999999996
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 524.6667
$ws.Range("I126").Value = 400.9
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 1202.7
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = 1267.3
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 35388.848
$ws.Range("I132").Value = 31246.94
$ws.Range("J132").Value = 40645.883
$ws.Range("K132").Value = 93740.81999999999
$ws.Range("L132").Value = 121937.649
$ws.Range("M132").Value = -91210.81999999999
$ws.Range("N132").Value = -126997.649

$ws.Range("H136").Value = 37298.86
$ws.Range("I136").Value = 24461.047
$ws.Range("K136").Value = 73383.141
$ws.Range("M136").Value = -70833.141
